# Update FTSE 100 ticker list:
#  - Remove "British Land" (BLND) row, shifting subsequent rows up by one.
#  - Insert a new row for "IMI" (Engineering) right after "IHG" (alphabetical spot),
#    shifting the rows that follow back down by one.
# Net effect: rows 21-46 keep the data that used to sit one row below them,
# row 47 becomes the new IMI entry, and rows 48 onward are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("BT-A", "BT Group", "Fixed Line Telecommunications"),
    @("BNZL", "Bunzl", "Support Services"),
    @("BRBY", "Burberry", "Personal Goods"),
    @("CNA", "Centrica", "Utilities"),
    @("CCH", "Coca-Cola HBC", "Beverages"),
    @("CPG", "Compass Group", "Support Services"),
    @("CTEC", "Convatec", "Health Care"),
    @("CRH", "CRH plc", "Construction & Materials"),
    @("CRDA", "Croda International", "Chemicals"),
    @("DCC", "DCC plc", "Support Services"),
    @("DGE", "Diageo", "Beverages"),
    @("EDV", "Endeavour Mining", "Mining"),
    @("ENT", "Entain", "Travel & Leisure"),
    @("EXPN", "Experian", "Support Services"),
    @("FCIT", "Foreign & Colonial Investment Trust", "Financial Services"),
    @("FLTR", "Flutter Entertainment", "Travel & Leisure"),
    @("FRAS", "Frasers Group", "Retail"),
    @("FRES", "Fresnillo plc", "Mining"),
    @("GLEN", "Glencore", "Mining"),
    @("GSK", "GSK plc", "Pharmaceuticals & Biotechnology"),
    @("HLN", "Haleon", "Pharmaceuticals & Biotechnology"),
    @("HLMA", "Halma plc", "Electronic & Electrical Equipment"),
    @("HL", "Hargreaves Lansdown", "Financial Services"),
    @("HSX", "Hiscox", "Insurance"),
    @("HSBA", "HSBC", "Banks"),
    @("IHG", "IHG Hotels & Resorts", "Travel & Leisure"),
    @("IMI", "IMI", "Engineering")
)

$startRow = 21
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
